$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindNewCarTest")

$ws.Range("A5").Value = "chrome"
$ws.Range("B5").Value = "tata"
$ws.Range("C5").Value = "Tata Cars"

$ws.Range("A6").Value = "firefox"
$ws.Range("B6").Value = "maruti"
$ws.Range("C6").Value = "Maruti Cars"

$ws.Range("A7").Value = "chrome"
$ws.Range("B7").Value = "honda"
$ws.Range("C7").Value = "Honda Cars"

$ws.Activate()
$ws.Range("D7").Select()
